# The slide has a single slide containing a TextBox ("TextBox 13", id=14)
# that currently sits right after the "Table 7" graphicFrame (shape index 2).
# The target OOXML moves that same shape (unchanged content) to become the
# very last shape in the slide's shape tree (after "TextBox 97").
#
# In the PowerPoint object model, moving a shape to be the last element of
# the spTree (i.e. drawn on top / at the end of z-order) is done by sending
# it to the front of the z-order stack with msoBringToFront (0).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $s.Shapes.Item("TextBox 13")
$shape.ZOrder(0)
